$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Belice" / "Nueva Caledonia" rows (192 = Nueva Caledonia, 193 = Belice)
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("A193").Value = "Belice"

$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Update "Datos actualizados..." timestamp text
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 10:04"

# Row 8 - Rusia
$ws.Range("B8").Value = 209688
$ws.Range("C8").Value = 11012
$ws.Range("D8").Value = 34306
$ws.Range("E8").Value = 173467
$ws.Range("G8").Value = 88
$ws.Range("H8").Value = 1915

# Row 88 - Eslovaquia
$ws.Range("B88").Value = 1457
$ws.Range("C88").Value = 2
$ws.Range("D88").Value = 941
$ws.Range("E88").Value = 490

# Row 104
$ws.Range("D104").Value = 321
$ws.Range("E104").Value = 517
